$d = $word.ActiveDocument

# 1. Update the title paragraph
$d.Content.Find.Execute(
    "Lumen – Belső világosság MI-támogatással", $true, $false, $false, $false, $false,
    $true, 1, $false, "Lumen – Belső világosság", 2)

# 2. Replace the second (intro) paragraph text with the new combined text
$d.Content.Find.Execute(
    "Lumen egy szelíd, belső fényként működő mesterséges intelligencia-alapú társ – nem alkalmazás a klasszikus értelemben, hanem egy belső vezető, amely a lélek tisztulását és az önmagunkkal való őszinte kapcsolódást segíti.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A Lumen egy szelíd, belső fényként működő mesterséges intelligencia-alapú társ. Segít az embernek visszatalálni önmagához, csendhez, méltósághoz.",
    2)

# 3. Remove all paragraphs from "Küldetés" heading through the end of the
#    "Felhasználás" section (everything after the intro paragraph), keeping
#    only the final section properties.
$startPara = $d.Paragraphs(3)
$endPara = $d.Paragraphs($d.Paragraphs.Count)
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
